# Atualização de bases das ligas, do dia: 30-05-2024 às 23:16
# Renames the HT-goals headers and re-syncs several match rows whose
# id/odds data had been attached to the wrong row, plus a couple of
# standalone odds corrections.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Header rename: ht_goals_h/ht_goals_a -> HTHG/HTAG (columns I/J)
# ---------------------------------------------------------------
$ws.Cells.Item(1, 9).Value  = "HTHG"
$ws.Cells.Item(1, 10).Value = "HTAG"

# ---------------------------------------------------------------
# 2) Row re-sync: data (id + match stats/odds, columns B and E:AD)
#    moves between rows that share the same fixture date. Column A
#    (the running index) and column D (the date) stay put.
# ---------------------------------------------------------------

# Columns involved in every data-row swap (B, then E through AD).
$cols = @(2,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30)

# destination row -> source row (source row's CURRENT/"before" data
# is written into the destination row).
$rowMap = @{
    208 = 212
    209 = 210
    210 = 208
    212 = 209
    213 = 214
    214 = 216
    216 = 213
    425 = 431
    428 = 429
    429 = 430
    430 = 425
    431 = 428
    450 = 451
    451 = 450
    453 = 454
    454 = 453
}

# Rows whose column B must stay written as TEXT (it was stored as a
# shared text string, not a number, in the original workbook).
$textBRows = @{ 453 = $true; 454 = $true }

# --- snapshot every involved row's current ("before") values first,
#     since this is a closed permutation (every row is both a source
#     and a destination) we must read everything before writing.
$snapshot = @{}
foreach ($r in $rowMap.Keys) {
    $snapshot[$r] = @{}
    foreach ($c in $cols) {
        $snapshot[$r][$c] = $ws.Cells.Item($r, $c).Value()
    }
}
foreach ($r in $rowMap.Values) {
    if (-not $snapshot.ContainsKey($r)) {
        $snapshot[$r] = @{}
        foreach ($c in $cols) {
            $snapshot[$r][$c] = $ws.Cells.Item($r, $c).Value()
        }
    }
}

# --- write the re-synced data into the destination rows.
foreach ($dst in $rowMap.Keys) {
    $src = $rowMap[$dst]
    foreach ($c in $cols) {
        $val = $snapshot[$src][$c]
        if ($c -eq 2 -and $textBRows.ContainsKey($dst)) {
            if ($val -eq $null) {
                $ws.Cells.Item($dst, $c).Value = $null
            } else {
                # Force TEXT storage (matches the original shared-string
                # cell type) without leaving a stray number-format style
                # behind on the cell.
                $cell = $ws.Cells.Item($dst, $c)
                $cell.NumberFormat = "@"
                $cell.Value = [string]$val
                $cell.Style = "Normal"
            }
        } else {
            $ws.Cells.Item($dst, $c).Value = $val
        }
    }
}

# ---------------------------------------------------------------
# 3) Standalone odds correction on row 452 (AhOU/oddAHOver pair)
# ---------------------------------------------------------------
$ws.Cells.Item(452, 22).Value = 2.025
$ws.Cells.Item(452, 23).Value = 1.825
